$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.904.31'
$ws.Range("E2").Value = '  -2.35%  '
$ws.Range("D3").Value = '1.756.66'
$ws.Range("E3").Value = '  -4.72%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9947'
$ws.Range("E4").Value = '  -0.63%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.19'
$ws.Range("E5").Value = '  -8.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9936'
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5007'
$ws.Range("E7").Value = '  -6.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.70'
$ws.Range("E8").Value = '  -7.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2643'
$ws.Range("E9").Value = '  -13.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06187'
$ws.Range("E10").Value = '  -10.22%  '
$ws.Range("D11").Value = '1.747.91'
$ws.Range("E11").Value = '  -5.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06920'
$ws.Range("E12").Value = '  -11.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.54'
$ws.Range("E13").Value = '  -15.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.490'
$ws.Range("E14").Value = '  -10.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5939'
$ws.Range("E15").Value = '  -21.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.30'
$ws.Range("E16").Value = '  -14.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9950'
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9963'
$ws.Range("E18").Value = '  -0.47%  '
$ws.Range("D19").Value = '25.910.65'
$ws.Range("E19").Value = '  -2.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.66'
$ws.Range("E20").Value = '  -16.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006759'
$ws.Range("E21").Value = '  -14.96%  '
$ws.Range("D22").Value = '1.967.90'
$ws.Range("E22").Value = '  -5.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.072'
$ws.Range("E23").Value = '  -11.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.073'
$ws.Range("E24").Value = '  -13.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.109'
$ws.Range("E25").Value = '  -14.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.62'
$ws.Range("E26").Value = '  -3.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.534'
$ws.Range("E27").Value = '  -8.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.847'
$ws.Range("E28").Value = '  -15.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.91'
$ws.Range("E29").Value = '  -12.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.57'
$ws.Range("E30").Value = '  -7.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.768'
$ws.Range("E31").Value = '  -11.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08066'
$ws.Range("E32").Value = '  -8.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.475'
$ws.Range("E33").Value = '  -15.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04480'
$ws.Range("E34").Value = '  -6.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9932'
$ws.Range("E35").Value = '  -0.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.616'
$ws.Range("E36").Value = '  -10.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9927'
$ws.Range("E37").Value = '  -12.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6056'
$ws.Range("E38").Value = '  -17.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.685'
$ws.Range("E39").Value = '  -13.64%  '
$ws.Range("E40").Value = '  -15.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '105.06'
$ws.Range("E41").Value = '  -3.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01523'
$ws.Range("E42").Value = '  -11.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9952'
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3814'
$ws.Range("E44").Value = '  -20.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.151'
$ws.Range("E45").Value = '  -12.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7335'
$ws.Range("E46").Value = '  -19.24%  '
$ws.Range("E47").Value = '  -10.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05220'
$ws.Range("E48").Value = '  -9.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.968'
$ws.Range("E50").Value = '  -13.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.51'
$ws.Range("E51").Value = '  -12.94%  '
